# Update code for report co so
# Adds two new sheets ("Đơn thu nợ" and "Lương") after the existing
# "Đơn sale chính" sheet, and fills "Đơn sale chính" with this month's
# sale-order data plus its totals row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Đơn sale chính" — fill header / data / totals rows
# ---------------------------------------------------------------
$wsSale = $wb.Worksheets.Item(1)
$wsSale.Name = "Đơn sale chính"

$saleHeaders = @(
    "Tiền tố", "Mã dịch vụ", "Ngày thực hiện", "Cơ sở", "Khách hàng",
    "Nguồn khách", "Tên dịch vụ", "Sale chính", "Đơn giá gốc", "Sale phụ",
    "Upsale", "Đơn giá", "Thanh toán lần đầu", "Trả sau", "Đã thanh toán",
    "Dư nợ", "Bác sĩ 1", "Bác sĩ 2", "Phụ phẫu 1", "Phụ phẫu 2",
    "Công phụ phẫu 1", "Công phụ phẫu 2"
)
for ($i = 0; $i -lt $saleHeaders.Length; $i++) {
    $wsSale.Cells.Item(1, $i + 1).Value = $saleHeaders[$i]
}

# Row 2 — single sale order for the period
$wsSale.Cells.Item(2, 1).Value = "HD-LUXURY"
$wsSale.Cells.Item(2, 2).Value = 529
$wsSale.Cells.Item(2, 3).NumberFormat = "@"
$wsSale.Cells.Item(2, 3).Value = "07-07-2024"
$wsSale.Cells.Item(2, 4).Value = "SÓC TRĂNG"
$wsSale.Cells.Item(2, 5).Value = "nguyễn thị hân"
$wsSale.Cells.Item(2, 6).Value = "Cá nhân"
$wsSale.Cells.Item(2, 7).Value = "Cắt mí"
$wsSale.Cells.Item(2, 8).Value = "Lê Đình Hậu"
$wsSale.Cells.Item(2, 9).Value = 7000000
$wsSale.Cells.Item(2, 12).Value = 7000000
$wsSale.Cells.Item(2, 13).Value = 4000000
$wsSale.Cells.Item(2, 14).Value = 0
$wsSale.Cells.Item(2, 15).Value = 4000000
$wsSale.Cells.Item(2, 16).Value = 3000000
$wsSale.Cells.Item(2, 17).Value = "Nguyễn Hoàng Yến Quyên"
$wsSale.Cells.Item(2, 19).Value = "Kha Như Huỳnh "

# Row 3 — totals
$wsSale.Cells.Item(3, 1).Value = "Tổng"
$wsSale.Cells.Item(3, 2).Value = 1
$wsSale.Cells.Item(3, 9).Value = 7000000
$wsSale.Cells.Item(3, 11).Value = 0
$wsSale.Cells.Item(3, 12).Value = 7000000
$wsSale.Cells.Item(3, 13).Value = 4000000
$wsSale.Cells.Item(3, 14).Value = 0
$wsSale.Cells.Item(3, 15).Value = 4000000
$wsSale.Cells.Item(3, 16).Value = 3000000
$wsSale.Cells.Item(3, 21).Value = 0
$wsSale.Cells.Item(3, 22).Value = 0

# ---------------------------------------------------------------
# Sheet 2: "Đơn thu nợ" — debt-collection orders
# ---------------------------------------------------------------
$wsDebt = $wb.Worksheets.Add($null, $wsSale)
$wsDebt.Name = "Đơn thu nợ"

$debtHeaders = @("Tiền tố", "Mã đơn thu nợ", "Đơn nợ", "Cơ sở", "Lượng thu", "Sale", "Ngày thu")
for ($i = 0; $i -lt $debtHeaders.Length; $i++) {
    $wsDebt.Cells.Item(1, $i + 1).Value = $debtHeaders[$i]
}

$wsDebt.Cells.Item(2, 1).Value = "TN"
$wsDebt.Cells.Item(2, 2).Value = 142
$wsDebt.Cells.Item(2, 3).Value = "HD-LUXURY-498"
$wsDebt.Cells.Item(2, 4).Value = "SÓC TRĂNG"
$wsDebt.Cells.Item(2, 5).Value = 2000000
$wsDebt.Cells.Item(2, 6).Value = "Lê Đình Hậu"
$wsDebt.Cells.Item(2, 7).NumberFormat = "@"
$wsDebt.Cells.Item(2, 7).Value = "07-07-2024"

$wsDebt.Cells.Item(3, 1).Value = "Tổng"
$wsDebt.Cells.Item(3, 2).Value = 1
$wsDebt.Cells.Item(3, 5).Value = 2000000

# ---------------------------------------------------------------
# Sheet 3: "Lương" — left blank, populated in a later update
# ---------------------------------------------------------------
$wsLuong = $wb.Worksheets.Add($null, $wsDebt)
$wsLuong.Name = "Lương"

$wsSale.Select()
